$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $NewValue)
    $rng = $ws.Range($CellRef)
    $rng.Value = "'" + $NewValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "29.465.40"
Set-TextValue "E2" "  +0.83%  "
Set-TextValue "D3" "1.874.35"
Set-TextValue "E3" "  +0.79%  "
Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  -0.01%  "
Set-TextValue "D5" "0.7175"
Set-TextValue "E5" "  +0.95%  "
Set-TextValue "D6" "239.37"
Set-TextValue "E6" "  +0.62%  "
Set-TextValue "D7" "1.0000"
Set-TextValue "E7" "  -0.09%  "
Set-TextValue "D8" "0.07840"
Set-TextValue "E8" "  -3.49%  "
Set-TextValue "D9" "0.3077"
Set-TextValue "E9" "  +1.16%  "
Set-TextValue "D10" "25.35"
Set-TextValue "E10" "  +9.38%  "
Set-TextValue "D11" "0.08248"
Set-TextValue "E11" "  +0.75%  "
Set-TextValue "D12" "1.885.36"
Set-TextValue "E12" "  +0.87%  "
Set-TextValue "D13" "5.245"
Set-TextValue "E13" "  +1.40%  "
Set-TextValue "D14" "0.7231"
Set-TextValue "E14" "  +2.33%  "
Set-TextValue "D15" "90.70"
Set-TextValue "E15" "  +1.27%  "
Set-TextValue "D16" "29.513.53"
Set-TextValue "E16" "  +0.91%  "
Set-TextValue "D17" "5.852"
Set-TextValue "E17" "  +1.07%  "
Set-TextValue "D18" "0.000007862"
Set-TextValue "E18" "  -0.42%  "
Set-TextValue "D19" "242.11"
Set-TextValue "E19" "  +2.12%  "
Set-TextValue "D20" "13.28"
Set-TextValue "E20" "  -0.45%  "
Set-TextValue "D21" "2.130.56"
Set-TextValue "E21" "  +0.21%  "
Set-TextValue "D22" "0.9992"
Set-TextValue "E22" "  -0.16%  "
Set-TextValue "E23" "  -0.02%  "
Set-TextValue "D24" "7.754"
Set-TextValue "E24" "  +4.40%  "
Set-TextValue "D25" "0.1558"
Set-TextValue "E25" "  +6.61%  "
Set-TextValue "D26" "163.13"
Set-TextValue "E26" "  +0.37%  "
Set-TextValue "D27" "9.018"
Set-TextValue "E27" "  +0.57%  "
Set-TextValue "D28" "18.39"
Set-TextValue "E28" "  +1.61%  "
Set-TextValue "D29" "1.937"
Set-TextValue "E29" "  -1.03%  "
Set-TextValue "D30" "1.356"
Set-TextValue "E30" "  -5.10%  "
Set-TextValue "E31" "  +0.01%  "
Set-TextValue "D32" "4.343"
Set-TextValue "E32" "  -1.26%  "
Set-TextValue "D33" "4.094"
Set-TextValue "E33" "  +1.97%  "
Set-TextValue "D34" "0.05264"
Set-TextValue "E34" "  +0.83%  "
Set-TextValue "D35" "1.201"
Set-TextValue "E35" "  +2.83%  "
Set-TextValue "D36" "0.7192"
Set-TextValue "E36" "  +1.65%  "
Set-TextValue "D37" "1.004"
Set-TextValue "E37" "  +0.31%  "
Set-TextValue "D38" "2.674"
Set-TextValue "E38" "  -0.04%  "
Set-TextValue "E39" "  +0.47%  "
Set-TextValue "D40" "2.722"
Set-TextValue "E40" "  -0.24%  "
Set-TextValue "D41" "1.178.37"
Set-TextValue "E41" "  +3.29%  "
Set-TextValue "D42" "0.9080"
Set-TextValue "E42" "  -1.56%  "
Set-TextValue "B43" "Aave"
Set-TextValue "C43" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D43" "72.36"
Set-TextValue "E43" "  +3.06%  "
Set-TextValue "B44" "FraxShare"
Set-TextValue "C44" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D44" "6.016"
Set-TextValue "E44" "  +2.45%  "
Set-TextValue "D45" "0.4316"
Set-TextValue "E45" "  +0.80%  "
Set-TextValue "D46" "0.9998"
Set-TextValue "E46" "  -0.05%  "
Set-TextValue "D47" "102.46"
Set-TextValue "E47" "  -0.05%  "
Set-TextValue "D48" "0.5366"
Set-TextValue "E48" "  -0.71%  "
Set-TextValue "D49" "1.768"
Set-TextValue "E49" "  -0.44%  "
Set-TextValue "D50" "9.170"
Set-TextValue "E50" "  -0.42%  "
Set-TextValue "D51" "7.033"
Set-TextValue "E51" "  +1.18%  "
